$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginCredentials")

# Capture the existing "Hyperlink" cell style before touching anything so it
# can be re-applied after the hyperlinks are rebuilt below.
$hyperlinkStyle = $ws1.Range("B2").Style

$rngB2 = $ws1.Range("B2")
$rngB3 = $ws1.Range("B3")

# Re-create the hyperlinks on B2/B3: B2 keeps its original mailto target
# (Jul@2020) but now needs an explicit display caption since the cell text
# itself is about to change; B3 is simply restored as-is.
$rngB2.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($rngB2, "mailto:Jul@2020", [Type]::Missing, [Type]::Missing, "Jul@2020")
$ws1.Hyperlinks.Add($rngB3, "mailto:July@123", [Type]::Missing, [Type]::Missing, [Type]::Missing)

$rngB2.Style = $hyperlinkStyle
$rngB3.Style = $hyperlinkStyle

# --- Update LoginCredentials rows 2 & 3 (agent id, password, BA/ALL flag) ---
$ws1.Range("A2").Value = 2388192
$ws1.Range("C2").Value = "BA"

$ws1.Range("A3").Value = 2390495
$ws1.Range("C3").Value = "ALL"

# Update B3 first, then B2, so the shared-string table ends up with
# "July@123$" before "secure#321" (matches original authoring order).
$ws1.Range("B3").Value = "July@123$"
$ws1.Range("B2").Value = "secure#321"

# --- Make LoginCredentials the active/selected sheet (was Priority) ---
$ws1.Activate()
